# Update the "Survey 4" row (row 4) of data in Sheet1 and move the
# active-cell selection, matching the re-uploaded version of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 ("Survey 4") responses changed: Very Happy, Quite Happy, Neutral,
# Quite Unhappy, Very Unhappy
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 40
$ws.Range("D4").Value = 7
$ws.Range("E4").Value = 7
$ws.Range("F4").Value = 5

# The sheet's active cell/selection moved from D7 to C8
$ws.Range("C8").Select()
